# Apply updates described by the commit:
#  - "Đơn phụ phẫu 1" sheet: insert two new order rows (HD-LUXURY 620 / 621)
#    before the "Tổng" (total) row, and update the total row figures.
#  - "Lương" sheet: recompute the base-salary / phụ phẫu 1 / total salary
#    figures for CẦN THƠ and HỆ THỐNG to reflect the new orders.

$wb = $excel.ActiveWorkbook

# ---- Sheet: "Đơn phụ phẫu 1" ----
$wsPhuPhau = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Move the existing "Tổng" row (currently row 3) down to row 5 and insert
# the two new data rows in its place.
$wsPhuPhau.Rows.Item(3).Insert()
$wsPhuPhau.Rows.Item(3).Insert()

$wsPhuPhau.Range("A3").Value = "HD-LUXURY"
$wsPhuPhau.Range("B3").Value = 620
$wsPhuPhau.Range("D3").Value = "CẦN THƠ"
$wsPhuPhau.Range("E3").Value = "Trần Thị Thanh Nhàn"
$wsPhuPhau.Range("F3").Value = "Cá nhân"
$wsPhuPhau.Range("G3").Value = "Nâng mũi"
$wsPhuPhau.Range("H3").Value = "Lâm Hoàng Phú"
$wsPhuPhau.Range("I3").Value = 100000

$wsPhuPhau.Range("A4").Value = "HD-LUXURY"
$wsPhuPhau.Range("B4").Value = 621
$wsPhuPhau.Range("D4").Value = "CẦN THƠ"
$wsPhuPhau.Range("E4").Value = "Trần Thị Ngọc Dung"
$wsPhuPhau.Range("F4").Value = "Cá nhân"
$wsPhuPhau.Range("G4").Value = "Nâng mũi"
$wsPhuPhau.Range("H4").Value = "Lâm Hoàng Phú"
$wsPhuPhau.Range("I4").Value = 100000

# The "Ngày thực hiện" column holds plain text dates (e.g. "08-01-2024" in
# row 2), not real date serials. Writing a date-shaped string through
# Range.Value auto-converts it to a date, so force text format first, then
# re-copy row 2's plain (unstyled) format on top so the new cells don't end
# up carrying a leftover text-number-format style.
$wsPhuPhau.Range("C3:C4").NumberFormat = "@"
$wsPhuPhau.Range("C3").Value = "08-03-2024"
$wsPhuPhau.Range("C4").Value = "08-03-2024"
$wsPhuPhau.Range("C2").Copy()
$wsPhuPhau.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the (now shifted-down) totals row 5. The insert above already
# carried the "Tổng" row's originally-blank C3:H3 cells down to C5:H5, so
# only the label/count/sum cells need new values here.
$wsPhuPhau.Range("B5").Value = 3
$wsPhuPhau.Range("I5").Value = 250000

# ---- Sheet: "Lương" ----
$wsLuong = $wb.Worksheets.Item("Lương")

$wsLuong.Range("B3").Value = 232142.8571428571
$wsLuong.Range("B8").Value = 250000
$wsLuong.Range("B31").Value = 482142.8571428572
$wsLuong.Range("B34").Value = 482142.8571428572
